$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column map: A=Loss, B=Depth, C=Architecture, D=Pooling, E=Type, F=Feature_num, G=specify, J=row index helper

# Row 186: only the J helper column continues
$ws.Cells.Item(186, 10).Value = 186

# Row 187: only the J helper column continues
$ws.Cells.Item(187, 10).Value = 187

# Rows 188-194 share the same Loss/Depth/Architecture/Pooling/Type/Feature_num combination
# (BCE, 2, GCN, gm, same, 12) with an increasing "specify" value in column G
$specifyValues = @(1, 2, 3, 4, 5, 6, 25)
for ($k = 0; $k -lt $specifyValues.Length; $k++) {
    $r = 188 + $k
    $ws.Cells.Item($r, 1).Value = "BCE"
    $ws.Cells.Item($r, 2).Value = 2
    $ws.Cells.Item($r, 3).Value = "GCN"
    $ws.Cells.Item($r, 4).Value = "gm"
    $ws.Cells.Item($r, 5).Value = "same"
    $ws.Cells.Item($r, 6).Value = 12
    $ws.Cells.Item($r, 7).Value = $specifyValues[$k]
    $ws.Cells.Item($r, 10).Value = $r
}

# Row 195: same combination but Pooling switches to "max" and no "specify" value
$ws.Cells.Item(195, 1).Value = "BCE"
$ws.Cells.Item(195, 2).Value = 2
$ws.Cells.Item(195, 3).Value = "GCN"
$ws.Cells.Item(195, 4).Value = "max"
$ws.Cells.Item(195, 5).Value = "same"
$ws.Cells.Item(195, 6).Value = 12
$ws.Cells.Item(195, 10).Value = 195

# Update the visible selection to match the end state of the edit session
$null = $ws.Range("G194").Select()

Write-Host "Applied noise combination rows 186-195"
